# Applies the Kraken_Profits market-data refresh across all leve sheets.
# Values were produced by the scheduled data-refresh runner; this script
# just re-pokes the updated numbers (and clears now-empty cells) into the
# corresponding worksheet/row/column cells.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30000
$ws.Range("J3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30228
$ws.Range("H6").Value = 31.25
$ws.Range("I6").Value = 31.25
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 93.75
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 18.25
$ws.Range("N6").Value = $null
$ws.Range("H32").Value = 9070.714
$ws.Range("J32").Value = 9070.714
$ws.Range("L32").Value = 9070.714
$ws.Range("N32").Value = -9722.714
$ws.Range("H33").Value = 358.57895
$ws.Range("I33").Value = 100.611115
$ws.Range("K33").Value = 100.611115
$ws.Range("M33").Value = 128.388885
$ws.Range("H58").Value = 200
$ws.Range("J58").Value = 200
$ws.Range("L58").Value = 600
$ws.Range("N58").Value = -900
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -36490
$ws.Range("H137").Value = 1624.8334
$ws.Range("I137").Value = 1450
$ws.Range("J137").Value = 2499
$ws.Range("K137").Value = 4350
$ws.Range("L137").Value = 7497
$ws.Range("M137").Value = -1800
$ws.Range("N137").Value = -12597
$ws.Range("H138").Value = 3826
$ws.Range("J138").Value = 3955.6365
$ws.Range("L138").Value = 11866.9095
$ws.Range("N138").Value = -22146.9095

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 500000
$ws.Range("I39").Value = 500000
$ws.Range("K39").Value = 500000
$ws.Range("M39").Value = -499480
$ws.Range("H45").Value = 3633
$ws.Range("I45").Value = 3633
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3633
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3256
$ws.Range("N45").Value = $null
$ws.Range("H74").Value = 1959.8
$ws.Range("I74").Value = 1959.8
$ws.Range("K74").Value = 1959.8
$ws.Range("M74").Value = -1085.8
$ws.Range("H77").Value = 1959.8
$ws.Range("I77").Value = 1959.8
$ws.Range("K77").Value = 9799
$ws.Range("M77").Value = -5431

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -875
$ws.Range("N30").Value = $null
$ws.Range("H107").Value = 234
$ws.Range("I107").Value = 234
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 234
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1686
$ws.Range("N107").Value = $null
$ws.Range("H108").Value = 99988
$ws.Range("J108").Value = 99988
$ws.Range("L108").Value = 99988
$ws.Range("N108").Value = -107668

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 2760
$ws.Range("I45").Value = 2760
$ws.Range("K45").Value = 2760
$ws.Range("M45").Value = -2167
$ws.Range("H60").Value = 19000
$ws.Range("I60").Value = 16000
$ws.Range("K60").Value = 16000
$ws.Range("M60").Value = -15489
$ws.Range("H95").Value = 25833.334
$ws.Range("J95").Value = 25833.334
$ws.Range("L95").Value = 25833.334
$ws.Range("N95").Value = -31325.334
$ws.Range("H98").Value = 99995
$ws.Range("J98").Value = 99995
$ws.Range("L98").Value = 99995
$ws.Range("N98").Value = -104487
$ws.Range("H99").Value = 1333.3334
$ws.Range("I99").Value = 1333.3334
$ws.Range("K99").Value = 1333.3334
$ws.Range("M99").Value = 164.6666
$ws.Range("H105").Value = 1642.7142
$ws.Range("I105").Value = 562.25
$ws.Range("K105").Value = 562.25
$ws.Range("M105").Value = 1184.75
$ws.Range("H126").Value = 1333.3334
$ws.Range("I126").Value = 1333.3334
$ws.Range("K126").Value = 4000.0002
$ws.Range("M126").Value = -1530.0002

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1072.7
$ws.Range("I34").Value = 685.4
$ws.Range("J34").Value = 1460
$ws.Range("K34").Value = 2056.2
$ws.Range("L34").Value = 4380
$ws.Range("M34").Value = -1972.2
$ws.Range("N34").Value = -4548
$ws.Range("H39").Value = 4557
$ws.Range("I39").Value = 3392.3333
$ws.Range("J39").Value = 4945.222
$ws.Range("K39").Value = 10176.9999
$ws.Range("L39").Value = 14835.666
$ws.Range("M39").Value = -9882.999899999999
$ws.Range("N39").Value = -15423.666
$ws.Range("H45").Value = 1033
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1033
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3099
$ws.Range("M45").Value = $null
$ws.Range("N45").Value = -4163
$ws.Range("H55").Value = 3526.875
$ws.Range("J55").Value = 3745
$ws.Range("L55").Value = 11235
$ws.Range("N55").Value = -11589
$ws.Range("H82").Value = 3999.5
$ws.Range("J82").Value = 5000
$ws.Range("L82").Value = 15000
$ws.Range("N82").Value = -15812
$ws.Range("H85").Value = 3999.5
$ws.Range("J85").Value = 5000
$ws.Range("L85").Value = 15000
$ws.Range("N85").Value = -17808
$ws.Range("H129").Value = 2527.9285
$ws.Range("I129").Value = 1355.8572
$ws.Range("J129").Value = 3700
$ws.Range("K129").Value = 4067.5716
$ws.Range("L129").Value = 11100
$ws.Range("M129").Value = 932.4284000000002
$ws.Range("N129").Value = -21100
$ws.Range("H131").Value = 1017.61536
$ws.Range("I131").Value = 653
$ws.Range("K131").Value = 1959
$ws.Range("M131").Value = 3081
$ws.Range("H138").Value = 3015
$ws.Range("I138").Value = 3016.6667
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 9050.000100000001
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -3910.000100000001
$ws.Range("N138").Value = -19280

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H44").Value = 5000000
$ws.Range("I44").Value = 5000000
$ws.Range("K44").Value = 5000000
$ws.Range("M44").Value = -4999404
$ws.Range("H49").Value = 35000
$ws.Range("J49").Value = 35000
$ws.Range("L49").Value = 35000
$ws.Range("N49").Value = -35368
$ws.Range("H70").Value = 3500
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 2000
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -1730
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 3500
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -1064
$ws.Range("N73").Value = -6872
$ws.Range("H80").Value = 69235
$ws.Range("I80").Value = 3849.5
$ws.Range("J80").Value = 200006
$ws.Range("K80").Value = 3849.5
$ws.Range("L80").Value = 200006
$ws.Range("M80").Value = -2851.5
$ws.Range("N80").Value = -202002
$ws.Range("H83").Value = 69235
$ws.Range("I83").Value = 3849.5
$ws.Range("J83").Value = 200006
$ws.Range("K83").Value = 19247.5
$ws.Range("L83").Value = 1000030
$ws.Range("M83").Value = -14255.5
$ws.Range("N83").Value = -1010014
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 5500
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 16500
$ws.Range("M132").Value = -13970

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5622.75
$ws.Range("I7").Value = 5622.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5622.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5510.75
$ws.Range("N7").Value = $null
$ws.Range("H46").Value = 3224.0645
$ws.Range("I46").Value = 1725.5
$ws.Range("J46").Value = 3327.4138
$ws.Range("K46").Value = 1725.5
$ws.Range("L46").Value = 3327.4138
$ws.Range("M46").Value = -1537.5
$ws.Range("N46").Value = -3703.4138
$ws.Range("H126").Value = 5622.75
$ws.Range("I126").Value = 5622.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16868.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14398.25
$ws.Range("N126").Value = $null
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5333.3335
$ws.Range("K136").Value = 16000.0005
$ws.Range("M136").Value = -13450.0005

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4061.2222
$ws.Range("I126").Value = 2936
$ws.Range("K126").Value = 8808
$ws.Range("M126").Value = -6338
